$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 11835.692
$ws.Range("J107").Value = 10000
$ws.Range("L107").Value = 10000
$ws.Range("N107").Value = -13840

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9718.5
$ws.Range("I45").Value = 11933
$ws.Range("K45").Value = 11933
$ws.Range("M45").Value = -11556

$ws.Range("H61").Value = 3782.56
$ws.Range("J61").Value = 4816.0625
$ws.Range("L61").Value = 4816.0625
$ws.Range("N61").Value = -5240.0625

$ws.Range("H74").Value = 2675.8438
$ws.Range("I74").Value = 1692.8235
$ws.Range("K74").Value = 1692.8235
$ws.Range("M74").Value = -818.8235

$ws.Range("H77").Value = 2675.8438
$ws.Range("I77").Value = 1692.8235
$ws.Range("K77").Value = 8464.1175
$ws.Range("M77").Value = -4096.1175

$ws.Range("H102").Value = 9539692
$ws.Range("I102").Value = 15896987
$ws.Range("K102").Value = 15896987
$ws.Range("M102").Value = -15895365

$ws.Range("H122").Value = 1226851.8
$ws.Range("I122").Value = 2787.5557
$ws.Range("J122").Value = 2450916
$ws.Range("K122").Value = 8362.667099999999
$ws.Range("L122").Value = 7352748
$ws.Range("M122").Value = -5912.667099999999
$ws.Range("N122").Value = -7357648

$ws.Range("H136").Value = 3782.56
$ws.Range("J136").Value = 4816.0625
$ws.Range("L136").Value = 14448.1875
$ws.Range("N136").Value = -19548.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4393.2354
$ws.Range("I107").Value = 3907.0833
$ws.Range("K107").Value = 3907.0833
$ws.Range("M107").Value = -1987.0833

$ws.Range("H122").Value = 44036.668
$ws.Range("J122").Value = 44036.668
$ws.Range("L122").Value = 44036.668
$ws.Range("N122").Value = -53836.668

$ws.Range("H134").Value = 6129.375
$ws.Range("I134").Value = 5008.4
$ws.Range("K134").Value = 15025.2
$ws.Range("M134").Value = -12490.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 10666
$ws.Range("I41").Value = 2000
$ws.Range("J41").Value = 14999
$ws.Range("K41").Value = 2000
$ws.Range("L41").Value = 14999
$ws.Range("M41").Value = -1572
$ws.Range("N41").Value = -15855

$ws.Range("H51").Value = 6500
$ws.Range("I51").Value = 6500
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 6500
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -5764
$ws.Range("N51").ClearContents()

$ws.Range("H59").Value = 74997
$ws.Range("J59").Value = 74997
$ws.Range("L59").Value = 74997
$ws.Range("N59").Value = -77287

$ws.Range("H60").Value = 3664.3333
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H61").Value = 6500
$ws.Range("I61").Value = 6500
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6500
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -6152
$ws.Range("N61").ClearContents()

$ws.Range("H132").Value = 17646.96
$ws.Range("I132").Value = 1242.1111
$ws.Range("K132").Value = 3726.3333
$ws.Range("M132").Value = -1196.3333

$ws.Range("H134").Value = 4185.5356
$ws.Range("I134").Value = 2200.3333
$ws.Range("J134").Value = 10141.143
$ws.Range("K134").Value = 6600.999899999999
$ws.Range("L134").Value = 30423.429
$ws.Range("M134").Value = -4065.999899999999
$ws.Range("N134").Value = -35493.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 24536606
$ws.Range("I4").Value = 3890286.2
$ws.Range("K4").Value = 11670858.6
$ws.Range("M4").Value = -11670746.6

$ws.Range("H110").Value = 33299.668
$ws.Range("I110").Value = 33299.668
$ws.Range("K110").Value = 99899.00399999999
$ws.Range("M110").Value = -95809.00399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 18965
$ws.Range("I122").Value = 18000
$ws.Range("J122").Value = 19930
$ws.Range("K122").Value = 54000
$ws.Range("L122").Value = 59790
$ws.Range("M122").Value = -51550
$ws.Range("N122").Value = -64690

$ws.Range("H132").Value = 4412.8184
$ws.Range("J132").Value = 4706
$ws.Range("L132").Value = 14118
$ws.Range("N132").Value = -19178

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 24918.957
$ws.Range("I7").Value = 37398.848
$ws.Range("K7").Value = 37398.848
$ws.Range("M7").Value = -37286.848

$ws.Range("H16").Value = 3368.5833
$ws.Range("I16").Value = 2878.3333
$ws.Range("K16").Value = 2878.3333
$ws.Range("M16").Value = -2708.3333

$ws.Range("H61").Value = 10332.444
$ws.Range("J61").Value = 9873.625
$ws.Range("L61").Value = 9873.625
$ws.Range("N61").Value = -10277.625

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()

$ws.Range("H113").Value = 10332.444
$ws.Range("J113").Value = 9873.625
$ws.Range("L113").Value = 9873.625
$ws.Range("N113").Value = -14213.625

$ws.Range("H122").Value = 6121.385
$ws.Range("I122").Value = 4332
$ws.Range("J122").Value = 6658.2
$ws.Range("K122").Value = 12996
$ws.Range("L122").Value = 19974.6
$ws.Range("M122").Value = -10546
$ws.Range("N122").Value = -24874.6

$ws.Range("H126").Value = 24918.957
$ws.Range("I126").Value = 37398.848
$ws.Range("K126").Value = 112196.544
$ws.Range("M126").Value = -109726.544

$ws.Range("H136").Value = 4430.4653
$ws.Range("I136").Value = 1625.4348
$ws.Range("J136").Value = 7656.25
$ws.Range("K136").Value = 4876.3044
$ws.Range("L136").Value = 22968.75
$ws.Range("M136").Value = -2326.3044
$ws.Range("N136").Value = -28068.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 18500
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 18500
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 18500
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -19644

$ws.Range("H100").Value = 24826.73
$ws.Range("I100").Value = 16763
$ws.Range("J100").Value = 46714
$ws.Range("K100").Value = 33526
$ws.Range("L100").Value = 93428
$ws.Range("M100").Value = -32985
$ws.Range("N100").Value = -94510

$ws.Range("H107").Value = 88821.28999999999
$ws.Range("I107").Value = 4350
$ws.Range("J107").Value = 299999.5
$ws.Range("K107").Value = 13050
$ws.Range("L107").Value = 899998.5
$ws.Range("M107").Value = -11130
$ws.Range("N107").Value = -903838.5

$ws.Range("H113").Value = 1447.9615
$ws.Range("I113").Value = 615.7368
$ws.Range("K113").Value = 1847.2104
$ws.Range("M113").Value = 322.7896000000001

$ws.Range("H122").Value = 13523.945
$ws.Range("I122").Value = 2371.054
$ws.Range("J122").Value = 36449.332
$ws.Range("K122").Value = 7113.162
$ws.Range("L122").Value = 109347.996
$ws.Range("M122").Value = -4663.162
$ws.Range("N122").Value = -114247.996
